# Automatic update of files.
# Adds a display-text second argument to the HYPERLINK() formulas in
# columns S, T, V, W, X, Y for rows 2-4 of the active sheet.
#
# Note: the diff being reproduced contains a pre-existing authoring bug:
# for columns T, V, W, X, Y the closing quote of the URL string was not
# added before the "; <label>" part was appended (only column S got a
# correctly closed string). We replicate that exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$base = "https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND"

# row number -> case id used both in the URL and as the hyperlink label
$rows = @{
    2 = "A 30683-2023"
    3 = "A 32699-2023"
    4 = "A 29992-2023"
}

# column -> (subfolder, file extension, whether the URL quote is closed
# before the appended "; label" text)
$cols = @(
    @{ Col = "S"; Folder = "artfynd";         Ext = "xlsx"; Closed = $true  },
    @{ Col = "T"; Folder = "kartor";          Ext = "png";  Closed = $false },
    @{ Col = "V"; Folder = "klagomål";        Ext = "docx"; Closed = $false },
    @{ Col = "W"; Folder = "klagomålsmail";   Ext = "docx"; Closed = $false },
    @{ Col = "X"; Folder = "tillsyn";         Ext = "docx"; Closed = $false },
    @{ Col = "Y"; Folder = "tillsynsmail";    Ext = "docx"; Closed = $false }
)

foreach ($rowNum in $rows.Keys) {
    $caseId = $rows[$rowNum]
    foreach ($c in $cols) {
        $url = "$base/$($c.Folder)/$caseId.$($c.Ext)"
        if ($c.Closed) {
            $formula = '=HYPERLINK("' + $url + '"; "' + $caseId + '")'
        } else {
            $formula = '=HYPERLINK("' + $url + '; "' + $caseId + '")'
        }
        $cellRef = "$($c.Col)$rowNum"
        $ws.Range($cellRef).Formula = $formula
    }
}
